$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apr 5 data came in, so row 39 (4/5/2020) stops being a forecast row and
# becomes an "actual" row. First, match the cell's look-and-feel to the
# other actual-data cells (e.g. I38) by copying over just its formatting
# (fill color / number format) rather than clobbering the whole cell.
$ws.Range("I38").Copy()
$ws.Range("I39").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# The forecast chain (I39:I48 in the old workbook) now has to start one row
# later, off of the new actual value in I39, so rebuild the forecast
# formulas for I40:I48 anchored there.
$ws.Range("I40:I48").Formula = "=I39*(1+AVERAGE(M37:M39))"

# Now drop in the actual reported US confirmed-case count for 4/5/2020.
$ws.Range("I39").Value = 336673

# The user's selection ends up one row down from before, on the newly
# vacated first forecast cell.
$ws.Range("I40").Select()
